# Adds a "FineshedAt" timestamp column to each of the five data sheets,
# mirroring a test run that wrote a Company/Contact/Country/whole-table
# column and then stamped when the write finished.

$wb = $excel.ActiveWorkbook

# Matches the quantised column-width grid the engine snaps custom widths to;
# the closest achievable value to the recorded "FineshedAt" bestFit width.
$finishedAtColWidth = 10.3

function Add-FinishedAtColumn($SheetName, $HeaderCol, $TimestampValue) {
    $ws = $wb.Worksheets.Item($SheetName)
    $headerCell = $ws.Cells.Item(1, $HeaderCol)
    $valueCell = $ws.Cells.Item(4, $HeaderCol)

    $headerCell.Value = "FineshedAt"
    $valueCell.Value = $TimestampValue

    # Give the new header cell the same shaded header look the other
    # header cells on the row already use (fillId=3 / indexed color 55).
    $headerCell.Interior.ColorIndex = 48
    $headerCell.Interior.Pattern = 1

    $ws.Columns.Item($HeaderCol).ColumnWidth = $finishedAtColWidth

    # The timestamp text embeds a newline; undo the engine's automatic
    # row-height bump so row 4 keeps its default (unflagged) height.
    $ws.Rows.Item(4).AutoFit()
}

# Sheet "writeCompanyColumnIntoXcel": Company column in A, new column in B.
Add-FinishedAtColumn "writeCompanyColumnIntoXcel" 2 "1571334563820`nThu Oct 17 10:49:23 PDT 2019"

# Sheet "writeContactColumnIntoXcel": Contact column in A, new column in B.
Add-FinishedAtColumn "writeContactColumnIntoXcel" 2 "1571334565652`nThu Oct 17 10:49:25 PDT 2019"

# Sheet "writeCountryColumnIntoXcel": Country column in A, new column in B.
Add-FinishedAtColumn "writeCountryColumnIntoXcel" 2 "1571334566494`nThu Oct 17 10:49:26 PDT 2019"

# Sheet "writeWholeTableNestedFor": Company/Contact/Country in A/B/C, new column in D.
Add-FinishedAtColumn "writeWholeTableNestedFor" 4 "1571334569471`nThu Oct 17 10:49:29 PDT 2019"

# Sheet "writeWholeTableSingleFor": Company/Contact/Country in A/B/C, new column in D.
Add-FinishedAtColumn "writeWholeTableSingleFor" 4 "1571334570770`nThu Oct 17 10:49:30 PDT 2019"
